# Auto-generated edit script for Pandaemonium_Profits (FFXIV leve profit tracker)
# Applies numeric corrections to columns H-N across ALC, ARM, BSM, CRP, CUL, GSM, WVR sheets
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(64, 8).Value = 411984.97  # H64: was 411986.3
$ws.Cells.Item(64, 9).Value = 788249.3  # I64: was 788250.75
$ws.Cells.Item(64, 10).Value = 4365.25  # J64: was 4366.5
$ws.Cells.Item(64, 11).Value = 788249.3  # K64: was 788250.75
$ws.Cells.Item(64, 12).Value = 4365.25  # L64: was 4366.5
$ws.Cells.Item(64, 13).Value = -788001.3  # M64: was -788002.75
$ws.Cells.Item(64, 14).Value = -4861.25  # N64: was -4862.5
$ws.Cells.Item(67, 8).Value = 411984.97  # H67: was 411986.3
$ws.Cells.Item(67, 9).Value = 788249.3  # I67: was 788250.75
$ws.Cells.Item(67, 10).Value = 4365.25  # J67: was 4366.5
$ws.Cells.Item(67, 11).Value = 788249.3  # K67: was 788250.75
$ws.Cells.Item(67, 12).Value = 4365.25  # L67: was 4366.5
$ws.Cells.Item(67, 13).Value = -787391.3  # M67: was -787392.75
$ws.Cells.Item(67, 14).Value = -6081.25  # N67: was -6082.5
$ws.Cells.Item(81, 8).Value = 38994.668  # H81: was 40328
$ws.Cells.Item(81, 10).Value = 38994.668  # J81: was 40328
$ws.Cells.Item(81, 12).Value = 38994.668  # L81: was 40328
$ws.Cells.Item(81, 14).Value = -40990.668  # N81: was -42324
$ws.Cells.Item(84, 8).Value = 38994.668  # H84: was 40328
$ws.Cells.Item(84, 10).Value = 38994.668  # J84: was 40328
$ws.Cells.Item(84, 12).Value = 116984.004  # L84: was 120984
$ws.Cells.Item(84, 14).Value = -126968.004  # N84: was -130968
$ws.Cells.Item(99, 8).Value = 857.2222  # H99: was 781.75
$ws.Cells.Item(99, 9).Value = 406.8  # I99: was 392.33334
$ws.Cells.Item(99, 10).Value = 1420.25  # J99: was 1950
$ws.Cells.Item(99, 11).Value = 1220.4  # K99: was 1177.00002
$ws.Cells.Item(99, 12).Value = 4260.75  # L99: was 5850
$ws.Cells.Item(99, 13).Value = 277.5999999999999  # M99: was 320.9999800000001
$ws.Cells.Item(99, 14).Value = -7256.75  # N99: was -8846
$ws.Cells.Item(101, 8).Value = 1725.091  # H101: was 1797.6
$ws.Cells.Item(101, 10).Value = 2405.7144  # J101: was 2640
$ws.Cells.Item(101, 12).Value = 7217.1432  # L101: was 7920
$ws.Cells.Item(101, 14).Value = -10461.1432  # N101: was -11164
$ws.Cells.Item(107, 8).Value = 516.4  # H107: was 521.93335
$ws.Cells.Item(107, 9).Value = 529.3333  # I107: was 535.4815
$ws.Cells.Item(107, 11).Value = 529.3333  # K107: was 535.4815
$ws.Cells.Item(107, 13).Value = 1390.6667  # M107: was 1384.5185
$ws.Cells.Item(113, 8).Value = 2345.2222  # H113: was 2443.423
$ws.Cells.Item(113, 9).Value = 1792.5294  # I113: was 1871.6666
$ws.Cells.Item(113, 10).Value = 3284.8  # J113: was 3223.0908
$ws.Cells.Item(113, 11).Value = 1792.5294  # K113: was 1871.6666
$ws.Cells.Item(113, 12).Value = 3284.8  # L113: was 3223.0908
$ws.Cells.Item(113, 13).Value = 1461.4706  # M113: was 1382.3334
$ws.Cells.Item(113, 14).Value = -9792.799999999999  # N113: was -9731.0908
$ws.Cells.Item(116, 8).Value = 2318.75  # H116: was 2305.8823
$ws.Cells.Item(116, 9).Value = 1983.3334  # I116: was 2000
$ws.Cells.Item(116, 10).Value = 2520  # J116: was 2577.7778
$ws.Cells.Item(116, 11).Value = 1983.3334  # K116: was 2000
$ws.Cells.Item(116, 12).Value = 2520  # L116: was 2577.7778
$ws.Cells.Item(116, 13).Value = 1458.6666  # M116: was 1442
$ws.Cells.Item(116, 14).Value = -9404  # N116: was -9461.7778
$ws.Cells.Item(132, 8).Value = 1318.0408  # H132: was 1144.8644
$ws.Cells.Item(132, 9).Value = 1341.9166  # I132: was 1161.638
$ws.Cells.Item(132, 11).Value = 4025.7498  # K132: was 3484.914
$ws.Cells.Item(132, 13).Value = -1495.7498  # M132: was -954.9139999999998
$ws.Cells.Item(138, 8).Value = 5112.822  # H138: was 5181.0137
$ws.Cells.Item(138, 9).Value = 1313.7931  # I138: was 1350.7142
$ws.Cells.Item(138, 10).Value = 7616.727  # J138: was 7618.477
$ws.Cells.Item(138, 11).Value = 3941.379300000001  # K138: was 4052.1426
$ws.Cells.Item(138, 12).Value = 22850.181  # L138: was 22855.431
$ws.Cells.Item(138, 13).Value = 1198.620699999999  # M138: was 1087.8574
$ws.Cells.Item(138, 14).Value = -33130.181  # N138: was -33135.431

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 19955.107  # H32: was 22180.932
$ws.Cells.Item(32, 9).Value = 22296.941  # I32: was 24729.13
$ws.Cells.Item(32, 10).Value = 11424.143  # J32: was 12412.833
$ws.Cells.Item(32, 11).Value = 22296.941  # K32: was 24729.13
$ws.Cells.Item(32, 12).Value = 11424.143  # L32: was 12412.833
$ws.Cells.Item(32, 13).Value = -22009.941  # M32: was -24442.13
$ws.Cells.Item(32, 14).Value = -11998.143  # N32: was -12986.833
$ws.Cells.Item(61, 8).Value = 5165.982  # H61: was 5740.12
$ws.Cells.Item(61, 9).Value = 3317.45  # I61: was 3754.3713
$ws.Cells.Item(61, 10).Value = 9787.3125  # J61: was 10373.533
$ws.Cells.Item(61, 11).Value = 3317.45  # K61: was 3754.3713
$ws.Cells.Item(61, 12).Value = 9787.3125  # L61: was 10373.533
$ws.Cells.Item(61, 13).Value = -3105.45  # M61: was -3542.3713
$ws.Cells.Item(61, 14).Value = -10211.3125  # N61: was -10797.533
$ws.Cells.Item(132, 8).Value = 2431.3076  # H132: was 2504.6
$ws.Cells.Item(132, 9).Value = 1632.6666  # I132: was 1693.4706
$ws.Cells.Item(132, 11).Value = 4897.9998  # K132: was 5080.4118
$ws.Cells.Item(132, 13).Value = -2367.9998  # M132: was -2550.4118
$ws.Cells.Item(136, 8).Value = 5165.982  # H136: was 5740.12
$ws.Cells.Item(136, 9).Value = 3317.45  # I136: was 3754.3713
$ws.Cells.Item(136, 10).Value = 9787.3125  # J136: was 10373.533
$ws.Cells.Item(136, 11).Value = 9952.349999999999  # K136: was 11263.1139
$ws.Cells.Item(136, 12).Value = 29361.9375  # L136: was 31120.599
$ws.Cells.Item(136, 13).Value = -7402.349999999999  # M136: was -8713.1139
$ws.Cells.Item(136, 14).Value = -34461.9375  # N136: was -36220.599

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(62, 8).Value = 28000  # H62: was 36181
$ws.Cells.Item(62, 10).Value = 28000  # J62: was 36181
$ws.Cells.Item(62, 12).Value = 28000  # L62: was 36181
$ws.Cells.Item(62, 14).Value = -29372  # N62: was -37553
$ws.Cells.Item(63, 8).Value = 38271  # H63: was 40271
$ws.Cells.Item(63, 10).Value = 38271  # J63: was 40271
$ws.Cells.Item(63, 12).Value = 38271  # L63: was 40271
$ws.Cells.Item(63, 14).Value = -39643  # N63: was -41643
$ws.Cells.Item(65, 8).Value = 28000  # H65: was 36181
$ws.Cells.Item(65, 10).Value = 28000  # J65: was 36181
$ws.Cells.Item(65, 12).Value = 84000  # L65: was 108543
$ws.Cells.Item(65, 14).Value = -90864  # N65: was -115407
$ws.Cells.Item(66, 8).Value = 38271  # H66: was 40271
$ws.Cells.Item(66, 10).Value = 38271  # J66: was 40271
$ws.Cells.Item(66, 12).Value = 114813  # L66: was 120813
$ws.Cells.Item(66, 14).Value = -121677  # N66: was -127677
$ws.Cells.Item(68, 8).Value = 26765  # H68: was 40295
$ws.Cells.Item(68, 9).Value = 10000  # I68: was 0
$ws.Cells.Item(68, 10).Value = 35147.5  # J68: was 40295
$ws.Cells.Item(68, 11).Value = 10000  # K68: was 0
$ws.Cells.Item(68, 12).Value = 35147.5  # L68: was 40295
$ws.Cells.Item(68, 13).Value = -9189  # M68: was None
$ws.Cells.Item(68, 14).Value = -36769.5  # N68: was -41917
$ws.Cells.Item(69, 8).Value = 30000  # H69: was 36295
$ws.Cells.Item(69, 10).Value = 30000  # J69: was 36295
$ws.Cells.Item(69, 12).Value = 30000  # L69: was 36295
$ws.Cells.Item(69, 14).Value = -31622  # N69: was -37917
$ws.Cells.Item(71, 8).Value = 26765  # H71: was 40295
$ws.Cells.Item(71, 9).Value = 10000  # I71: was 0
$ws.Cells.Item(71, 10).Value = 35147.5  # J71: was 40295
$ws.Cells.Item(71, 11).Value = 30000  # K71: was 0
$ws.Cells.Item(71, 12).Value = 105442.5  # L71: was 120885
$ws.Cells.Item(71, 13).Value = -25944  # M71: was None
$ws.Cells.Item(71, 14).Value = -113554.5  # N71: was -128997
$ws.Cells.Item(72, 8).Value = 30000  # H72: was 36295
$ws.Cells.Item(72, 10).Value = 30000  # J72: was 36295
$ws.Cells.Item(72, 12).Value = 90000  # L72: was 108885
$ws.Cells.Item(72, 14).Value = -98112  # N72: was -116997
$ws.Cells.Item(75, 8).Value = 21408.334  # H75: was 13000
$ws.Cells.Item(75, 9).Value = 14214  # I75: was 10000
$ws.Cells.Item(75, 10).Value = 22847.2  # J75: was 13750
$ws.Cells.Item(75, 11).Value = 14214  # K75: was 10000
$ws.Cells.Item(75, 12).Value = 22847.2  # L75: was 13750
$ws.Cells.Item(75, 13).Value = -13278  # M75: was -9064
$ws.Cells.Item(75, 14).Value = -24719.2  # N75: was -15622
$ws.Cells.Item(76, 8).Value = 34588.4  # H76: was 34788.4
$ws.Cells.Item(76, 10).Value = 34588.4  # J76: was 34788.4
$ws.Cells.Item(76, 12).Value = 34588.4  # L76: was 34788.4
$ws.Cells.Item(76, 14).Value = -35218.4  # N76: was -35418.4
$ws.Cells.Item(78, 8).Value = 21408.334  # H78: was 13000
$ws.Cells.Item(78, 9).Value = 14214  # I78: was 10000
$ws.Cells.Item(78, 10).Value = 22847.2  # J78: was 13750
$ws.Cells.Item(78, 11).Value = 42642  # K78: was 30000
$ws.Cells.Item(78, 12).Value = 68541.60000000001  # L78: was 41250
$ws.Cells.Item(78, 13).Value = -37962  # M78: was -25320
$ws.Cells.Item(78, 14).Value = -77901.60000000001  # N78: was -50610
$ws.Cells.Item(79, 8).Value = 34588.4  # H79: was 34788.4
$ws.Cells.Item(79, 10).Value = 34588.4  # J79: was 34788.4
$ws.Cells.Item(79, 12).Value = 34588.4  # L79: was 34788.4
$ws.Cells.Item(79, 14).Value = -36772.4  # N79: was -36972.4
$ws.Cells.Item(82, 8).Value = 23587  # H82: was 11691.375
$ws.Cells.Item(82, 9).Value = 2291.6  # I82: was 2227.5
$ws.Cells.Item(82, 10).Value = 38798  # J82: was 40083
$ws.Cells.Item(82, 11).Value = 2291.6  # K82: was 2227.5
$ws.Cells.Item(82, 12).Value = 38798  # L82: was 40083
$ws.Cells.Item(82, 13).Value = -1908.6  # M82: was -1844.5
$ws.Cells.Item(82, 14).Value = -39564  # N82: was -40849
$ws.Cells.Item(85, 8).Value = 23587  # H85: was 11691.375
$ws.Cells.Item(85, 9).Value = 2291.6  # I85: was 2227.5
$ws.Cells.Item(85, 10).Value = 38798  # J85: was 40083
$ws.Cells.Item(85, 11).Value = 2291.6  # K85: was 2227.5
$ws.Cells.Item(85, 12).Value = 38798  # L85: was 40083
$ws.Cells.Item(85, 13).Value = -965.5999999999999  # M85: was -901.5
$ws.Cells.Item(85, 14).Value = -41450  # N85: was -42735
$ws.Cells.Item(94, 8).Value = 1566.6666  # H94: was 823.13043
$ws.Cells.Item(94, 9).Value = 1566.6666  # I94: was 849.4666999999999
$ws.Cells.Item(94, 10).Value = 0  # J94: was 773.75
$ws.Cells.Item(94, 11).Value = 1566.6666  # K94: was 849.4666999999999
$ws.Cells.Item(94, 12).Value = 0  # L94: was 773.75
$ws.Cells.Item(94, 13).Value = -1115.6666  # M94: was -398.4666999999999
$ws.Cells.Item(94, 14).ClearContents()  # N94: was -1675.75
$ws.Cells.Item(134, 8).Value = 4753.9287  # H134: was 3496.7058
$ws.Cells.Item(134, 9).Value = 5232.273  # I134: was 3818.7693
$ws.Cells.Item(134, 10).Value = 3000  # J134: was 2450
$ws.Cells.Item(134, 11).Value = 15696.819  # K134: was 11456.3079
$ws.Cells.Item(134, 12).Value = 9000  # L134: was 7350
$ws.Cells.Item(134, 13).Value = -13161.819  # M134: was -8921.3079
$ws.Cells.Item(134, 14).Value = -14070  # N134: was -12420

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 4419.5347  # H31: was 4688.85
$ws.Cells.Item(31, 9).Value = 5018.846  # I31: was 5212.16
$ws.Cells.Item(31, 10).Value = 3502.9412  # J31: was 3816.6667
$ws.Cells.Item(31, 11).Value = 5018.846  # K31: was 5212.16
$ws.Cells.Item(31, 12).Value = 3502.9412  # L31: was 3816.6667
$ws.Cells.Item(31, 13).Value = -4723.846  # M31: was -4917.16
$ws.Cells.Item(31, 14).Value = -4092.9412  # N31: was -4406.6667
$ws.Cells.Item(34, 8).Value = 4419.5347  # H34: was 4688.85
$ws.Cells.Item(34, 9).Value = 5018.846  # I34: was 5212.16
$ws.Cells.Item(34, 10).Value = 3502.9412  # J34: was 3816.6667
$ws.Cells.Item(34, 11).Value = 5018.846  # K34: was 5212.16
$ws.Cells.Item(34, 12).Value = 3502.9412  # L34: was 3816.6667
$ws.Cells.Item(34, 13).Value = -4816.846  # M34: was -5010.16
$ws.Cells.Item(34, 14).Value = -3906.9412  # N34: was -4220.6667
$ws.Cells.Item(99, 8).Value = 1980.6364  # H99: was 2271.913
$ws.Cells.Item(99, 9).Value = 1475.8572  # I99: was 1777.8462
$ws.Cells.Item(99, 10).Value = 2864  # J99: was 2914.2
$ws.Cells.Item(99, 11).Value = 1475.8572  # K99: was 1777.8462
$ws.Cells.Item(99, 12).Value = 2864  # L99: was 2914.2
$ws.Cells.Item(99, 13).Value = 22.14280000000008  # M99: was -279.8462
$ws.Cells.Item(99, 14).Value = -5860  # N99: was -5910.2
$ws.Cells.Item(126, 8).Value = 1980.6364  # H126: was 2271.913
$ws.Cells.Item(126, 9).Value = 1475.8572  # I126: was 1777.8462
$ws.Cells.Item(126, 10).Value = 2864  # J126: was 2914.2
$ws.Cells.Item(126, 11).Value = 4427.571599999999  # K126: was 5333.5386
$ws.Cells.Item(126, 12).Value = 8592  # L126: was 8742.599999999999
$ws.Cells.Item(126, 13).Value = -1957.571599999999  # M126: was -2863.5386
$ws.Cells.Item(126, 14).Value = -13532  # N126: was -13682.6
$ws.Cells.Item(132, 8).Value = 1731.4889  # H132: was 1759.5454
$ws.Cells.Item(132, 9).Value = 1336.6389  # I132: was 1360.6285
$ws.Cells.Item(132, 11).Value = 4009.9167  # K132: was 4081.8855
$ws.Cells.Item(132, 13).Value = -1479.9167  # M132: was -1551.8855
$ws.Cells.Item(134, 8).Value = 3199.814  # H134: was 3136.182
$ws.Cells.Item(134, 9).Value = 1993.3334  # I134: was 1909.4736
$ws.Cells.Item(134, 11).Value = 5980.0002  # K134: was 5728.4208
$ws.Cells.Item(134, 13).Value = -3445.0002  # M134: was -3193.4208

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(97, 8).Value = 2107.8  # H97: was 2203.111
$ws.Cells.Item(97, 9).Value = 987.5  # I97: was 1300
$ws.Cells.Item(97, 10).Value = 2854.6667  # J97: was 2461.1428
$ws.Cells.Item(97, 11).Value = 2962.5  # K97: was 3900
$ws.Cells.Item(97, 12).Value = 8564.000100000001  # L97: was 7383.428400000001
$ws.Cells.Item(97, 13).Value = -2466.5  # M97: was -3404
$ws.Cells.Item(97, 14).Value = -9556.000100000001  # N97: was -8375.428400000001
$ws.Cells.Item(113, 8).Value = 683.92957  # H113: was 685.92755
$ws.Cells.Item(113, 9).Value = 697.63416  # I113: was 701.87177
$ws.Cells.Item(113, 11).Value = 2092.90248  # K113: was 2105.61531
$ws.Cells.Item(113, 13).Value = 77.09752000000026  # M113: was 64.38468999999986

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 4088.6924  # H102: was 2825.195
$ws.Cells.Item(102, 9).Value = 4059.2222  # I102: was 2502.818
$ws.Cells.Item(102, 11).Value = 4059.2222  # K102: was 2502.818
$ws.Cells.Item(102, 13).Value = -2437.2222  # M102: was -880.8180000000002
$ws.Cells.Item(132, 8).Value = 2733  # H132: was 2626.0454
$ws.Cells.Item(132, 9).Value = 2498.3333  # I132: was 2365.9375
$ws.Cells.Item(132, 11).Value = 7494.999899999999  # K132: was 7097.8125
$ws.Cells.Item(132, 13).Value = -4964.999899999999  # M132: was -4567.8125

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(95, 8).Value = 0  # H95: was 40000
$ws.Cells.Item(95, 10).Value = 0  # J95: was 40000
$ws.Cells.Item(95, 12).Value = 0  # L95: was 40000
$ws.Cells.Item(95, 14).ClearContents()  # N95: was -45492
$ws.Cells.Item(122, 8).Value = 2000.5  # H122: was 2059.353
$ws.Cells.Item(122, 9).Value = 1500.2858  # I122: was 1538.7693
$ws.Cells.Item(122, 11).Value = 4500.857400000001  # K122: was 4616.3079
$ws.Cells.Item(122, 13).Value = -2050.857400000001  # M122: was -2166.3079
$ws.Cells.Item(136, 8).Value = 9676.941000000001  # H136: was 7437.6
$ws.Cells.Item(136, 9).Value = 8376.352999999999  # I136: was 5649
$ws.Cells.Item(136, 10).Value = 10977.529  # J136: was 9885.157999999999
$ws.Cells.Item(136, 11).Value = 25129.059  # K136: was 16947
$ws.Cells.Item(136, 12).Value = 32932.587  # L136: was 29655.474
$ws.Cells.Item(136, 13).Value = -22579.059  # M136: was -14397
$ws.Cells.Item(136, 14).Value = -38032.587  # N136: was -34755.474

Write-Host "Applied all Pandaemonium_Profits corrections."